$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: replace the contents of a Range with a raw OOXML fragment (a
# sequence of block-level <w:p> elements) via Range.InsertXML, which (per
# this host) REPLACES the range's current contents with the supplied XML.
# InsertXML requires the full mc-package envelope, not a bare fragment.
#
# NOTE: named-parameter binding (-Range ... -InnerXml ...) is unreliable in
# this PowerShell host when the first arg is a COM object, so this is always
# called with positional arguments.
# ---------------------------------------------------------------------------
function Set-RangeOoxml($Range, [string]$InnerXml) {
    $pkg = '<?xml version="1.0" standalone="yes"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' + $InnerXml + '</w:body>' +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'

    $Range.InsertXML($pkg) | Out-Null
}

# ---------------------------------------------------------------------------
# 1. Add a new "App List view, do not load icon in onBindView() ..." Heading4
#    bullet right before the empty spacer Heading4 paragraph that precedes
#    the "Good to have" Heading2 (i.e. right after the
#    "Tracker to see, how many notifications have been silent on a day."
#    Heading4 item).
# ---------------------------------------------------------------------------
$anchorPara = $null
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    $txt = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($txt -eq "Tracker to see, how many notifications have been silent on a day.") {
        $anchorPara = $p
        break
    }
}

if ($anchorPara -ne $null) {
    $anchorPara.Range.InsertParagraphAfter()

    # The newly inserted (empty, Heading4-styled) paragraph now sits right
    # after $anchorPara; re-fetch it by its (now known) index.
    $newIndex = $anchorPara.Index + 1
    $newRange = $d.Paragraphs.Item($newIndex).Range

    $newParaXml = '<w:p><w:pPr><w:pStyle w:val="Heading4"/></w:pPr>' +
        '<w:r><w:t xml:space="preserve">App List view, do not load icon in </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>onBindView</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t>(), instead load all icons on activity start and pass to AppListAdapter</w:t></w:r>' +
        '</w:p>'

    Set-RangeOoxml $newRange $newParaXml
}

# ---------------------------------------------------------------------------
# 2. Remove the stray <w:lastRenderedPageBreak/> from the "Transaction
#    across multiple Daos" Heading3 paragraph.
# ---------------------------------------------------------------------------
$targetPara = $null
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    $txt = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($txt -eq "Transaction across multiple Daos") {
        $targetPara = $p
        break
    }
}

if ($targetPara -ne $null) {
    $fixedParaXml = '<w:p><w:pPr><w:pStyle w:val="Heading3"/></w:pPr>' +
        '<w:r><w:t xml:space="preserve">Transaction across multiple </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>Daos</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '</w:p>'

    Set-RangeOoxml $targetPara.Range $fixedParaXml
}
